$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the description text in column B (rows 2-8) ---
$ws.Range("B2").Value = "Some tables are ignored (not described) in Excel workbooks for version 1A. Some constraints are dropped."
$ws.Range("B3").Value = "Some columns are changed from mandatory to optional. Two tables renamed. Tickets and TicketLines. Some tables are dropped. More constraints are dropped."
$ws.Range("B4").Value = "Some tables are renamed (including Products_*). Some relationships are changed."
$ws.Range("B5").Value = "Some columns are dropped: SFlag, SiteGUId, Attribute*."
$ws.Range("B6").Value = "Some columns are renamed (and merged [name]). Some constraints are renamed"
$ws.Range("B7").Value = "Columns are retyped."
$ws.Range("B8").Value = "Some tables are added and/or renamed (e.g, Product_kits; with different name)."

# --- Make the header row bold (new bold font + style) ---
$ws.Range("A1:B1").Font.Bold = $true

# --- Row 2 gets the same "wrapped description" look as the other rows ---
$ws.Range("A2").VerticalAlignment = -4108
$ws.Range("B2").VerticalAlignment = -4108
$ws.Range("B2").WrapText = $true
$ws.Rows(2).RowHeight = 28.8

# --- Widen column B to fit the longer text ---
$ws.Columns(2).ColumnWidth = 66.5

# --- Update the active selection ---
$ws.Range("B5").Select()
